$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 18, pushing existing rows 18:125 down to 19:126.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new data record.
$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(18, 3).Value = "La Araucanía"
$ws.Cells.Item(18, 4).Value = 44602
$ws.Cells.Item(18, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 9
$ws.Cells.Item(18, 6).Value = 100112012
$ws.Cells.Item(18, 7).Value = "Espinaca"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 55
$ws.Cells.Item(18, 11).Value = 13000
$ws.Cells.Item(18, 12).Value = 13000
$ws.Cells.Item(18, 13).Value = 13000
$ws.Cells.Item(18, 14).Value = "$/docena de atados"
$ws.Cells.Item(18, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(18, 16).Value = 4333
$ws.Cells.Item(18, 17).Value = 3
$ws.Cells.Item(18, 18).Value = "Hortaliza"
